# Apply updated crypto price/volume figures to Sheet1 (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain decimal number need to be forced to
# the Text format first, otherwise Excel auto-converts them to a float and the
# original formatted digits (e.g. trailing zero) get lost to FP rounding noise.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = '27.954.34'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '1.635.75'
$ws.Range("E3").Value = '  -1.08%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '212.28'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").Value = '0.523'
$ws.Range("E6").Value = '  -1.19%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '23.33'
$ws.Range("E8").Value = '  -1.18%  '
$ws.Range("D9").Value = '0.258'
$ws.Range("E9").Value = '  -2.76%  '
$ws.Range("E10").Value = '  -0.08%  '
$ws.Range("E11").Value = '  +1.03%  '
$ws.Range("D12").Value = '1.868.01'
$ws.Range("D13").Value = '1.636.05'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  -0.72%  '
$ws.Range("D15").Value = '0.565'
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("D16").Value = '65.44'
$ws.Range("E16").Value = '  -0.64%  '
$ws.Range("D17").Value = '27.961.20'
$ws.Range("E17").Value = '  -0.50%  '
$ws.Range("D18").Value = '231.05'
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("E19").Value = '  -0.24%  '
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("E22").Value = '  -1.03%  '
$ws.Range("E23").Value = '  -2.85%  '
$ws.Range("E24").Value = '  -3.84%  '
$ws.Range("D25").Value = '154.87'
$ws.Range("E25").Value = '  +1.53%  '
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("E28").Value = '  -1.18%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("E31").Value = '  -0.54%  '
$ws.Range("E32").Value = '  +1.37%  '
$ws.Range("D33").Value = '1.407.90'
$ws.Range("E33").Value = '  -3.16%  '
$ws.Range("E34").Value = '  -0.51%  '
$ws.Range("E35").Value = '  -0.55%  '
$ws.Range("D36").Value = '1.02'
$ws.Range("E36").Value = '  +9.41%  '
$ws.Range("E37").Value = '  +1.44%  '
$ws.Range("E38").Value = '  +0.43%  '
$ws.Range("D39").Value = '0.563'
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("D40").Value = '0.873'
$ws.Range("E40").Value = '  -2.36%  '
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '66.90'
$ws.Range("E43").Value = '  -3.70%  '
$ws.Range("D44").Value = '5.54'
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("D45").Value = '1.83'
$ws.Range("E45").Value = '  +0.08%  '
$ws.Range("E46").Value = '  -1.34%  '
$ws.Range("D47").Value = '1.777.15'
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("D48").Value = '87.93'
$ws.Range("E48").Value = '  -1.42%  '
$ws.Range("D49").Value = '0.0₆0106'
$ws.Range("E49").Value = '  +0.49%  '
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("E51").Value = '  -0.42%  '
